# Auto-generated edit script: updates cached market-price / profit values
# in the Goblin_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values come from a scheduled market-data refresh; no formulas are involved,
# all target cells hold plain cached numbers.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 10: A Jawbreaking Weapon of Staggering Weight | Whispering Maple Wand
$ws.Range("H10").Value = 5668
$ws.Range("I10").Value = 6002
$ws.Range("K10").Value = 6002
$ws.Range("M10").Value = -5709

# Row 58: A Matter of Vital Importance | Mega-Potion of Vitality
$ws.Range("H58").Value = 5655
$ws.Range("J58").Value = 7248.7144
$ws.Range("L58").Value = 21746.1432
$ws.Range("N58").Value = -22046.1432


$ws = $wb.Worksheets.Item("ARM")

# Row 6: Don't Hit Me One More Time | Bronze Hoplon
$ws.Range("H6").Value = 90169.22
$ws.Range("I6").Value = 181250
$ws.Range("K6").Value = 181250
$ws.Range("M6").Value = -181077

# Row 95: Shielded Life | High Steel Scutum
$ws.Range("H95").Value = 76502.8
$ws.Range("J95").Value = 76502.8
$ws.Range("L95").Value = 76502.8
$ws.Range("N95").Value = -81994.8

# Row 97: Ore for Me | High Steel Ingot
$ws.Range("H97").Value = 409.5
$ws.Range("I97").Value = 409.5
$ws.Range("K97").Value = 409.5
$ws.Range("M97").Value = 86.5

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 5921.304
$ws.Range("I132").Value = 6904.737
$ws.Range("J132").Value = 1250
$ws.Range("K132").Value = 20714.211
$ws.Range("L132").Value = 3750
$ws.Range("M132").Value = -18184.211
$ws.Range("N132").Value = -8810


$ws = $wb.Worksheets.Item("BSM")

# Row 7: Thank You for Your Business | Bronze Bastard Sword
$ws.Range("H7").Value = 1400
$ws.Range("I7").Value = 500
$ws.Range("K7").Value = 500
$ws.Range("M7").Value = -387

# Row 20: Smelt and Dealt | Iron Ingot
$ws.Range("H20").Value = 2571.6667
$ws.Range("I20").Value = 2703.0667
$ws.Range("J20").Value = 2243.1667
$ws.Range("K20").Value = 2703.0667
$ws.Range("L20").Value = 2243.1667
$ws.Range("M20").Value = -2456.0667
$ws.Range("N20").Value = -2737.1667

# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("I86").Value = 925.80646
$ws.Range("K86").Value = 925.80646
$ws.Range("M86").Value = 197.19354

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("I89").Value = 925.80646
$ws.Range("K89").Value = 4629.0323
$ws.Range("M89").Value = 986.9677000000001

# Row 94: High Steal | High Steel Nugget
$ws.Range("H94").Value = 1666.2368
$ws.Range("I94").Value = 1415.2354
$ws.Range("K94").Value = 1415.2354
$ws.Range("M94").Value = -964.2354

# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Range("H99").Value = 3275.64
$ws.Range("I99").Value = 3025
$ws.Range("K99").Value = 3025
$ws.Range("M99").Value = -1527

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 482782.12
$ws.Range("J134").Value = 5001720
$ws.Range("L134").Value = 15005160
$ws.Range("N134").Value = -15010230


$ws = $wb.Worksheets.Item("CRP")

# Row 8: Bows for the Boys | Maple Longbow
$ws.Range("H8").Value = 499
$ws.Range("J8").Value = 499
$ws.Range("L8").Value = 499
$ws.Range("N8").Value = -779

# Row 10: Spears and Sorcery | Maple Crook
$ws.Range("H10").Value = 1480.8462
$ws.Range("I10").Value = 312.625
$ws.Range("J10").Value = 3350
$ws.Range("K10").Value = 312.625
$ws.Range("L10").Value = 3350
$ws.Range("M10").Value = -173.625
$ws.Range("N10").Value = -3628

# Row 16: Raise the Roof | Ash Lumber
$ws.Range("H16").Value = 2939.4
$ws.Range("I16").Value = 2224.8333
$ws.Range("J16").Value = 3415.7778
$ws.Range("K16").Value = 2224.8333
$ws.Range("L16").Value = 3415.7778
$ws.Range("M16").Value = -1937.8333
$ws.Range("N16").Value = -3989.7778

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 2993.8333
$ws.Range("I31").Value = 1032.5385
$ws.Range("J31").Value = 4493.647
$ws.Range("K31").Value = 1032.5385
$ws.Range("L31").Value = 4493.647
$ws.Range("M31").Value = -737.5385000000001
$ws.Range("N31").Value = -5083.647

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 2993.8333
$ws.Range("I34").Value = 1032.5385
$ws.Range("J34").Value = 4493.647
$ws.Range("K34").Value = 1032.5385
$ws.Range("L34").Value = 4493.647
$ws.Range("M34").Value = -830.5385000000001
$ws.Range("N34").Value = -4897.647

# Row 50: The Arsenal of Theocracy | Cobalt Halberd
$ws.Range("H50").Value = 42855.57
$ws.Range("J50").Value = 43331.668
$ws.Range("L50").Value = 43331.668
$ws.Range("N50").Value = -44581.668

# Row 86: Birch, Please | Birch Lumber
$ws.Range("H86").Value = 5765.3335
$ws.Range("I86").Value = 4589.4
$ws.Range("K86").Value = 4589.4
$ws.Range("M86").Value = -3466.4

# Row 89: Built This City on Blocks and Soul (L) | Birch Lumber
$ws.Range("H89").Value = 5765.3335
$ws.Range("I89").Value = 4589.4
$ws.Range("K89").Value = 22947
$ws.Range("M89").Value = -17331

# Row 113: Patient Patients | White Ash Lumber
$ws.Range("H113").Value = 2939.4
$ws.Range("I113").Value = 2224.8333
$ws.Range("J113").Value = 3415.7778
$ws.Range("K113").Value = 2224.8333
$ws.Range("L113").Value = 3415.7778
$ws.Range("M113").Value = -54.83329999999978
$ws.Range("N113").Value = -7755.7778

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 1741.1082
$ws.Range("I132").Value = 1018.2647
$ws.Range("K132").Value = 3054.7941
$ws.Range("M132").Value = -524.7941000000001

# Row 133: Yimepi's Country Charms | Ginseng Earrings
$ws.Range("H133").Value = 95000
$ws.Range("J133").Value = 95000
$ws.Range("L133").Value = 95000
$ws.Range("N133").Value = -100060


$ws = $wb.Worksheets.Item("CUL")

# Row 94: All You Can Stomach | Baklava
$ws.Range("H94").Value = 9616.571
$ws.Range("J94").Value = 11215.333
$ws.Range("L94").Value = 33645.999
$ws.Range("N94").Value = -34997.999

# Row 114: One Last Meal | Mushroom Saute
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").ClearContents()
$ws.Range("M114").ClearContents()
$ws.Range("N114").Value = 0

# Row 121: A Cookie for Your Troubles | Coffee Biscuit
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = 0

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 2472289.8
$ws.Range("J131").Value = 3337100.8
$ws.Range("L131").Value = 10011302.4
$ws.Range("N131").Value = -10021382.4


$ws = $wb.Worksheets.Item("GSM")

# Row 20: Brothers in Arms | Brass Wristlets of Crafting
$ws.Range("H20").Value = 16000
$ws.Range("J20").Value = 22000
$ws.Range("L20").Value = 22000
$ws.Range("N20").Value = -22490

# Row 57: Gold Is So Last Year | Electrum Circlet (Amber)
$ws.Range("H57").Value = 28341.666

# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 40003748
$ws.Range("I80").Value = 58826264
$ws.Range("J80").Value = 5899.375
$ws.Range("K80").Value = 58826264
$ws.Range("L80").Value = 5899.375
$ws.Range("M80").Value = -58825266
$ws.Range("N80").Value = -7895.375

# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 40003748
$ws.Range("I83").Value = 58826264
$ws.Range("J83").Value = 5899.375
$ws.Range("K83").Value = 294131320
$ws.Range("L83").Value = 29496.875
$ws.Range("M83").Value = -294126328
$ws.Range("N83").Value = -39480.875

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 31259052
$ws.Range("J113").Value = 9978
$ws.Range("L113").Value = 9978
$ws.Range("N113").Value = -14318

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 52633290
$ws.Range("I132").Value = 62501844
$ws.Range("J132").Value = 999.6667
$ws.Range("K132").Value = 187505532
$ws.Range("L132").Value = 2999.0001
$ws.Range("M132").Value = -187503002
$ws.Range("N132").Value = -8059.0001


$ws = $wb.Worksheets.Item("LTW")

# Row 9: From the Sands to the Stage | Leather Himantes
$ws.Range("H9").Value = 245.41667
$ws.Range("J9").Value = 244.5
$ws.Range("L9").Value = 244.5
$ws.Range("N9").Value = -692.5

# Row 30: Packing a Punch | Goatskin Cesti
$ws.Range("H30").Value = 456.66666
$ws.Range("I30").Value = 456.66666
$ws.Range("K30").Value = 456.66666
$ws.Range("M30").Value = -348.66666

# Row 35: No Risk, No Reward | Toadskin Cesti
$ws.Range("H35").Value = 2302.2222
$ws.Range("I35").Value = 2302.2222
$ws.Range("K35").Value = 2302.2222
$ws.Range("M35").Value = -1966.2222

# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 4475.125
$ws.Range("I40").Value = 3246.5334
$ws.Range("K40").Value = 3246.5334
$ws.Range("M40").Value = -3110.5334

# Row 93: Hide to Go Seek | Gagana Leather
$ws.Range("H93").Value = 3693.532
$ws.Range("I93").Value = 1750.3043
$ws.Range("J93").Value = 5555.7915
$ws.Range("K93").Value = 1750.3043
$ws.Range("L93").Value = 5555.7915
$ws.Range("M93").Value = -502.3043
$ws.Range("N93").Value = -8051.7915

# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 12000
$ws.Range("N122").Value = -16900

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 2266.9412
$ws.Range("I132").Value = 2265.9697
$ws.Range("K132").Value = 6797.909100000001
$ws.Range("M132").Value = -4267.909100000001


$ws = $wb.Worksheets.Item("WVR")

# Row 14: Hat in Hand | Straw Hat
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = 0

# Row 62: Pride Up in Smoke | Rainbow Cloth
$ws.Range("H62").Value = 11375345
$ws.Range("J62").Value = 13776.441
$ws.Range("L62").Value = 13776.441
$ws.Range("N62").Value = -15024.441

# Row 65: Desperate for Diversionaries (L) | Rainbow Cloth
$ws.Range("H65").Value = 11375345
$ws.Range("J65").Value = 13776.441
$ws.Range("L65").Value = 68882.205
$ws.Range("N65").Value = -75122.205

# Row 70: An Account of My Boots | Holy Rainbow Shoes
$ws.Range("H70").Value = 500025000
$ws.Range("J70").Value = 500025000
$ws.Range("L70").Value = 500025000
$ws.Range("N70").Value = -500025630

# Row 73: Soot in My Hair and Scars on My Feet (L) | Holy Rainbow Shoes
$ws.Range("H73").Value = 500025000
$ws.Range("J73").Value = 500025000
$ws.Range("L73").Value = 500025000
$ws.Range("N73").Value = -500027184

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 15884751
$ws.Range("I132").Value = 16678614
$ws.Range("K132").Value = 50035842
$ws.Range("M132").Value = -50033312

